$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(5, 6).Value = 2
$ws.Cells.Item(15, 6).Value = 2
$ws.Cells.Item(16, 6).Value = 3
$ws.Cells.Item(17, 6).Value = -2
$ws.Cells.Item(18, 6).Value = -1
$ws.Cells.Item(20, 6).Value = 2
$ws.Cells.Item(28, 6).Value = -2
$ws.Cells.Item(30, 6).Value = -2
$ws.Cells.Item(35, 6).Value = -1
$ws.Cells.Item(41, 6).Value = 1
$ws.Cells.Item(44, 6).Value = 3
$ws.Cells.Item(46, 6).Value = -3
$ws.Cells.Item(51, 6).Value = 1
$ws.Cells.Item(56, 6).Value = -1
$ws.Cells.Item(61, 6).Value = -1
$ws.Cells.Item(67, 6).Value = -2
$ws.Cells.Item(71, 6).Value = -2
